$wb = $excel.ActiveWorkbook

# ALC row 12 (Leve Item ID 5515)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 408.92856
$ws.Range("I12").Value = 434.16666
$ws.Range("J12").Value = 257.5
$ws.Range("K12").Value = 434.16666
$ws.Range("L12").Value = 257.5
$ws.Range("M12").Value = -264.16666
$ws.Range("N12").Value = -597.5

# ALC row 28 (Leve Item ID 27772)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3513.111
$ws.Range("I28").Value = 295.14285
$ws.Range("J28").Value = 14776
$ws.Range("K28").Value = 295.14285
$ws.Range("L28").Value = 14776
$ws.Range("M28").Value = 189.85715
$ws.Range("N28").Value = -15746

# ALC row 121 (Leve Item ID 39731)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2929.75
$ws.Range("J121").Value = 2929.75
$ws.Range("L121").Value = 8789.25
$ws.Range("N121").Value = -12283.25

# ALC row 127 (Leve Item ID 36114)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2221.5
$ws.Range("I127").Value = 2221.5
$ws.Range("K127").Value = 6664.5
$ws.Range("M127").Value = -1704.5

# ALC row 129 (Leve Item ID 36115)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 3077.5715
$ws.Range("I129").Value = 3088.8333
$ws.Range("J129").Value = 3010
$ws.Range("K129").Value = 9266.499899999999
$ws.Range("L129").Value = 9030
$ws.Range("M129").Value = -4266.499899999999
$ws.Range("N129").Value = -19030

# ALC row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1925.18
$ws.Range("J138").Value = 1539.6316
$ws.Range("L138").Value = 4618.8948
$ws.Range("N138").Value = -14898.8948

# ARM row 28 (Leve Item ID 19534)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 55555
$ws.Range("I28").Value = 55555
$ws.Range("K28").Value = 55555
$ws.Range("M28").Value = -55363

# ARM row 41 (Leve Item ID 2501)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 18991.143
$ws.Range("I41").Value = 11019.667
$ws.Range("K41").Value = 11019.667
$ws.Range("M41").Value = -10605.667

# ARM row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4787.2
$ws.Range("I61").Value = 4456
$ws.Range("J61").Value = 4870
$ws.Range("K61").Value = 4456
$ws.Range("L61").Value = 4870
$ws.Range("M61").Value = -4244
$ws.Range("N61").Value = -5294

# ARM row 99 (Leve Item ID 19534)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 55555
$ws.Range("I99").Value = 55555
$ws.Range("K99").Value = 55555
$ws.Range("M99").Value = -52560

# ARM row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2316.8462
$ws.Range("I132").Value = 1547.0667
$ws.Range("K132").Value = 4641.2001
$ws.Range("M132").Value = -2111.2001

# ARM row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4787.2
$ws.Range("I136").Value = 4456
$ws.Range("J136").Value = 4870
$ws.Range("K136").Value = 13368
$ws.Range("L136").Value = 14610
$ws.Range("M136").Value = -10818
$ws.Range("N136").Value = -19710

# CRP row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4088.2856
$ws.Range("I31").Value = 2712.1667
$ws.Range("J31").Value = 12345
$ws.Range("K31").Value = 2712.1667
$ws.Range("L31").Value = 12345
$ws.Range("M31").Value = -2417.1667
$ws.Range("N31").Value = -12935

# CRP row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4088.2856
$ws.Range("I34").Value = 2712.1667
$ws.Range("J34").Value = 12345
$ws.Range("K34").Value = 2712.1667
$ws.Range("L34").Value = 12345
$ws.Range("M34").Value = -2510.1667
$ws.Range("N34").Value = -12749

# CRP row 52 (Leve Item ID 43237)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 66773
$ws.Range("J52").Value = 66773
$ws.Range("L52").Value = 66773
$ws.Range("N52").Value = -67361

# CRP row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3005
$ws.Range("I58").Value = 2781.4167
$ws.Range("J58").Value = 3899.3333
$ws.Range("K58").Value = 2781.4167
$ws.Range("L58").Value = 3899.3333
$ws.Range("M58").Value = -2578.4167
$ws.Range("N58").Value = -4305.3333

# CRP row 111 (Leve Item ID 25792)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 35000
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# CRP row 118 (Leve Item ID 26112)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 55000
$ws.Range("J118").Value = 55000
$ws.Range("L118").Value = 55000
$ws.Range("N118").Value = -58314

# CRP row 131 (Leve Item ID 35461)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 104999.5
$ws.Range("J131").Value = 104999.5
$ws.Range("L131").Value = 104999.5
$ws.Range("N131").Value = -115079.5

# CRP row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3052.9
$ws.Range("I134").Value = 3190.3333
$ws.Range("K134").Value = 9570.999899999999
$ws.Range("M134").Value = -7035.999899999999

# CRP row 136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3005
$ws.Range("I136").Value = 2781.4167
$ws.Range("J136").Value = 3899.3333
$ws.Range("K136").Value = 8344.250100000001
$ws.Range("L136").Value = 11697.9999
$ws.Range("M136").Value = -5794.250100000001
$ws.Range("N136").Value = -16797.9999

# CUL row 4 (Leve Item ID 4650)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 98115224
$ws.Range("I4").Value = 11862434
$ws.Range("J4").Value = 400000000
$ws.Range("K4").Value = 35587302
$ws.Range("L4").Value = 1200000000
$ws.Range("M4").Value = -35587190
$ws.Range("N4").Value = -1200000224

# CUL row 32 (Leve Item ID 4731)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 3000
$ws.Range("N32").Value = -3566

# CUL row 33 (Leve Item ID 4867)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 172.82353
$ws.Range("I33").Value = 173.4375
$ws.Range("J33").Value = 163
$ws.Range("K33").Value = 1040.625
$ws.Range("L33").Value = 978
$ws.Range("M33").Value = -757.625
$ws.Range("N33").Value = -1544

# CUL row 39 (Leve Item ID 4712)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 5079
$ws.Range("I39").Value = 748.5
$ws.Range("K39").Value = 2245.5
$ws.Range("M39").Value = -1951.5

# CUL row 51 (Leve Item ID 4646)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2144
$ws.Range("I51").Value = 2124
$ws.Range("K51").Value = 6372
$ws.Range("M51").Value = -5912

# CUL row 131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 41283.715
$ws.Range("I131").Value = 111949.3
$ws.Range("J131").Value = 2025.0555
$ws.Range("K131").Value = 335847.9
$ws.Range("L131").Value = 6075.166499999999
$ws.Range("M131").Value = -330807.9
$ws.Range("N131").Value = -16155.1665

# LTW row 22 (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 698.8077
$ws.Range("I22").Value = 1365.5
$ws.Range("J22").Value = 498.8
$ws.Range("K22").Value = 1365.5
$ws.Range("L22").Value = 498.8
$ws.Range("M22").Value = -1070.5
$ws.Range("N22").Value = -1088.8

# LTW row 27 (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 698.8077
$ws.Range("I27").Value = 1365.5
$ws.Range("J27").Value = 498.8
$ws.Range("K27").Value = 1365.5
$ws.Range("L27").Value = 498.8
$ws.Range("M27").Value = -1258.5
$ws.Range("N27").Value = -712.8

# LTW row 46 (Leve Item ID 5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2184.6843
$ws.Range("I46").Value = 1519.625
$ws.Range("J46").Value = 2668.3635
$ws.Range("K46").Value = 1519.625
$ws.Range("L46").Value = 2668.3635
$ws.Range("M46").Value = -1331.625
$ws.Range("N46").Value = -3044.3635

# LTW row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1730.4117
$ws.Range("I136").Value = 1067
$ws.Range("K136").Value = 3201
$ws.Range("M136").Value = -651

# WVR row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3278.0698
$ws.Range("I132").Value = 3500.5
$ws.Range("K132").Value = 10501.5
$ws.Range("M132").Value = -7971.5

# WVR row 136 (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 875.55554
$ws.Range("I136").Value = 875.55554
$ws.Range("K136").Value = 2626.66662
$ws.Range("M136").Value = -76.66661999999997
